$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 275
$ws.Range("I4").Value = 275
$ws.Range("K4").Value = 275
$ws.Range("M4").Value = -161
$ws.Range("H12").Value = 1601.3334
$ws.Range("I12").Value = 154.6
$ws.Range("K12").Value = 154.6
$ws.Range("M12").Value = 15.40000000000001
$ws.Range("H33").Value = 468
$ws.Range("I33").Value = 474.18182
$ws.Range("K33").Value = 474.18182
$ws.Range("M33").Value = -245.18182
$ws.Range("H132").Value = 5626.231
$ws.Range("I132").Value = 5664.6523
$ws.Range("J132").Value = 5331.6665
$ws.Range("K132").Value = 16993.9569
$ws.Range("L132").Value = 15994.9995
$ws.Range("M132").Value = -14463.9569
$ws.Range("N132").Value = -21054.9995
$ws.Range("H137").Value = 5138.857
$ws.Range("I137").Value = 3000
$ws.Range("J137").Value = 5303.385
$ws.Range("K137").Value = 9000
$ws.Range("L137").Value = 15910.155
$ws.Range("M137").Value = -6450
$ws.Range("N137").Value = -21010.155
$ws.Range("H138").Value = 7221.5
$ws.Range("J138").Value = 7581.9116
$ws.Range("L138").Value = 22745.7348
$ws.Range("N138").Value = -33025.73480000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H49").Value = 15000
$ws.Range("J49").Value = 15000
$ws.Range("L49").Value = 15000
$ws.Range("N49").Value = -15520
$ws.Range("H61").Value = 7203.0713
$ws.Range("I61").Value = 5908
$ws.Range("K61").Value = 5908
$ws.Range("M61").Value = -5696
$ws.Range("H63").Value = 9071.143
$ws.Range("I63").Value = 3499
$ws.Range("J63").Value = 9999.833000000001
$ws.Range("K63").Value = 3499
$ws.Range("L63").Value = 9999.833000000001
$ws.Range("M63").Value = -2813
$ws.Range("N63").Value = -11371.833
$ws.Range("H66").Value = 9071.143
$ws.Range("I66").Value = 3499
$ws.Range("J66").Value = 9999.833000000001
$ws.Range("K66").Value = 17495
$ws.Range("L66").Value = 49999.165
$ws.Range("M66").Value = -14063
$ws.Range("N66").Value = -56863.165
$ws.Range("H74").Value = 7333
$ws.Range("J74").Value = 9999.5
$ws.Range("L74").Value = 9999.5
$ws.Range("N74").Value = -11747.5
$ws.Range("H77").Value = 7333
$ws.Range("J77").Value = 9999.5
$ws.Range("L77").Value = 49997.5
$ws.Range("N77").Value = -58733.5
$ws.Range("H110").Value = 180853.42
$ws.Range("I110").Value = 201887.36
$ws.Range("K110").Value = 201887.36
$ws.Range("M110").Value = -199842.36
$ws.Range("H132").Value = 9957.950000000001
$ws.Range("I132").Value = 3554.9167
$ws.Range("J132").Value = 19562.5
$ws.Range("K132").Value = 10664.7501
$ws.Range("L132").Value = 58687.5
$ws.Range("M132").Value = -8134.750100000001
$ws.Range("N132").Value = -63747.5
$ws.Range("H136").Value = 7203.0713
$ws.Range("I136").Value = 5908
$ws.Range("K136").Value = 17724
$ws.Range("M136").Value = -15174
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1034.1111
$ws.Range("J64").Value = 1132.5
$ws.Range("L64").Value = 1132.5
$ws.Range("N64").Value = -1582.5
$ws.Range("H67").Value = 1034.1111
$ws.Range("J67").Value = 1132.5
$ws.Range("L67").Value = 1132.5
$ws.Range("N67").Value = -2692.5
$ws.Range("H80").Value = 1412.0526
$ws.Range("I80").Value = 1708.9
$ws.Range("J80").Value = 1082.2222
$ws.Range("K80").Value = 1708.9
$ws.Range("L80").Value = 1082.2222
$ws.Range("M80").Value = -710.9000000000001
$ws.Range("N80").Value = -3078.2222
$ws.Range("H83").Value = 1412.0526
$ws.Range("I83").Value = 1708.9
$ws.Range("J83").Value = 1082.2222
$ws.Range("K83").Value = 8544.5
$ws.Range("L83").Value = 5411.111
$ws.Range("M83").Value = -3552.5
$ws.Range("N83").Value = -15395.111
$ws.Range("H105").Value = 2436.2727
$ws.Range("I105").Value = 1900.3077
$ws.Range("J105").Value = 3210.4443
$ws.Range("K105").Value = 1900.3077
$ws.Range("L105").Value = 3210.4443
$ws.Range("M105").Value = -153.3077000000001
$ws.Range("N105").Value = -6704.4443
$ws.Range("H134").Value = 43009.926
$ws.Range("I134").Value = 5903.8335
$ws.Range("J134").Value = 117222.11
$ws.Range("K134").Value = 17711.5005
$ws.Range("L134").Value = 351666.33
$ws.Range("M134").Value = -15176.5005
$ws.Range("N134").Value = -356736.33
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 489
$ws.Range("J11").Value = 489
$ws.Range("L11").Value = 489
$ws.Range("N11").Value = -769
$ws.Range("H58").Value = 504910.75
$ws.Range("I58").Value = 1113802.5
$ws.Range("J58").Value = 6726.636
$ws.Range("K58").Value = 1113802.5
$ws.Range("L58").Value = 6726.636
$ws.Range("M58").Value = -1113599.5
$ws.Range("N58").Value = -7132.636
$ws.Range("H134").Value = 504075.34
$ws.Range("I134").Value = 3750.5
$ws.Range("K134").Value = 11251.5
$ws.Range("M134").Value = -8716.5
$ws.Range("H136").Value = 504910.75
$ws.Range("I136").Value = 1113802.5
$ws.Range("J136").Value = 6726.636
$ws.Range("K136").Value = 3341407.5
$ws.Range("L136").Value = 20179.908
$ws.Range("M136").Value = -3338857.5
$ws.Range("N136").Value = -25279.908
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 461.51852
$ws.Range("J26").Value = 491.64
$ws.Range("L26").Value = 1474.92
$ws.Range("N26").Value = -2050.92
$ws.Range("H56").Value = 6916.2
$ws.Range("I56").Value = 6916.2
$ws.Range("K56").Value = 6916.2
$ws.Range("M56").Value = -6386.2
$ws.Range("H131").Value = 4740.706
$ws.Range("J131").Value = 5983
$ws.Range("L131").Value = 17949
$ws.Range("N131").Value = -28029
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6819.643
$ws.Range("I70").Value = 6484.143
$ws.Range("K70").Value = 6484.143
$ws.Range("M70").Value = -6214.143
$ws.Range("H73").Value = 6819.643
$ws.Range("I73").Value = 6484.143
$ws.Range("K73").Value = 6484.143
$ws.Range("M73").Value = -5548.143
$ws.Range("H80").Value = 971673.9399999999
$ws.Range("I80").Value = 692392.9
$ws.Range("J80").Value = 1669876.6
$ws.Range("K80").Value = 692392.9
$ws.Range("L80").Value = 1669876.6
$ws.Range("M80").Value = -691394.9
$ws.Range("N80").Value = -1671872.6
$ws.Range("H83").Value = 971673.9399999999
$ws.Range("I83").Value = 692392.9
$ws.Range("J83").Value = 1669876.6
$ws.Range("K83").Value = 3461964.5
$ws.Range("L83").Value = 8349383
$ws.Range("M83").Value = -3456972.5
$ws.Range("N83").Value = -8359367
$ws.Range("H97").Value = 1536.5385
$ws.Range("I97").Value = 1515.5454
$ws.Range("K97").Value = 1515.5454
$ws.Range("M97").Value = -1019.5454
$ws.Range("H102").Value = 2879.56
$ws.Range("I102").Value = 2144.1177
$ws.Range("K102").Value = 2144.1177
$ws.Range("M102").Value = -522.1176999999998
$ws.Range("H132").Value = 691210.5
$ws.Range("J132").Value = 129362.875
$ws.Range("L132").Value = 388088.625
$ws.Range("N132").Value = -393148.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 630897.25
$ws.Range("I7").Value = 1004735.6
$ws.Range("K7").Value = 1004735.6
$ws.Range("M7").Value = -1004623.6
$ws.Range("H16").Value = 804.5833
$ws.Range("I16").Value = 757.7778
$ws.Range("J16").Value = 945
$ws.Range("K16").Value = 757.7778
$ws.Range("L16").Value = 945
$ws.Range("M16").Value = -587.7778
$ws.Range("N16").Value = -1285
$ws.Range("H22").Value = 964.06665
$ws.Range("I22").Value = 946.4
$ws.Range("K22").Value = 946.4
$ws.Range("M22").Value = -651.4
$ws.Range("H27").Value = 964.06665
$ws.Range("I27").Value = 946.4
$ws.Range("K27").Value = 946.4
$ws.Range("M27").Value = -839.4
$ws.Range("H40").Value = 719984.0600000001
$ws.Range("I40").Value = 1004377.8
$ws.Range("K40").Value = 1004377.8
$ws.Range("M40").Value = -1004241.8
$ws.Range("H46").Value = 3696.4517
$ws.Range("I46").Value = 3088.6843
$ws.Range("J46").Value = 4658.75
$ws.Range("K46").Value = 3088.6843
$ws.Range("L46").Value = 4658.75
$ws.Range("M46").Value = -2900.6843
$ws.Range("N46").Value = -5034.75
$ws.Range("H122").Value = 349073.44
$ws.Range("I122").Value = 4763.952
$ws.Range("K122").Value = 14291.856
$ws.Range("M122").Value = -11841.856
$ws.Range("H126").Value = 630897.25
$ws.Range("I126").Value = 1004735.6
$ws.Range("K126").Value = 3014206.8
$ws.Range("M126").Value = -3011736.8
$ws.Range("H132").Value = 7153.143
$ws.Range("I132").Value = 4838.8
$ws.Range("J132").Value = 8438.888999999999
$ws.Range("K132").Value = 14516.4
$ws.Range("L132").Value = 25316.667
$ws.Range("M132").Value = -11986.4
$ws.Range("N132").Value = -30376.667
$ws.Range("H139").Value = 55000
$ws.Range("J139").Value = 55000
$ws.Range("L139").Value = 55000
$ws.Range("N139").Value = -65280
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 2549.5
$ws.Range("J19").Value = 2549.5
$ws.Range("L19").Value = 2549.5
$ws.Range("N19").Value = -2897.5
$ws.Range("H126").Value = 3155.8572
$ws.Range("J126").Value = 5995.75
$ws.Range("L126").Value = 17987.25
$ws.Range("N126").Value = -22927.25
$ws.Range("H127").Value = 49666.668
$ws.Range("J127").Value = 49666.668
$ws.Range("L127").Value = 49666.668
$ws.Range("N127").Value = -59586.668
$ws.Range("H132").Value = 52457
$ws.Range("I132").Value = 3516.7273
$ws.Range("J132").Value = 106291.3
$ws.Range("K132").Value = 10550.1819
$ws.Range("L132").Value = 318873.9
$ws.Range("M132").Value = -8020.1819
$ws.Range("N132").Value = -323933.9
